$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.926.74"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.497.53"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").Value = "2.518.75"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "2.946.16"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "58.891.19"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "2.509.12"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "0.0₃0775"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  -1.72%  "
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.834"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "280.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.17%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
